$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update input data
$ws.Range("G12").Value = 230
$ws.Range("G13").Value = 85

# Update formulas that depend on the new data
$ws.Range("D19").Formula = "=G12*10+G13*10+G14"
$ws.Range("D20").Formula = "=G12*10+G13*10"
$ws.Range("D21").Formula = "=G12*10+G13*10"

$ws.Range("D22").Formula = "=G12*10+G13*10"
$ws.Range("C22").Formula = "=-D22*2.2"

$ws.Range("D23").Formula = "=G12*10+G13*10"
$ws.Range("C23").Formula = "=D23*2.2"

$ws.Range("D24").Formula = "=-4*D23"

# Update the active cell selection
$ws.Range("C24").Select()
